$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these "Price" cells to remain plain text (their values look like
# plain decimal numbers, e.g. "211.68", and Excel would otherwise silently
# convert them into floating point numbers when assigned via .Value).
$textPriceCells = @("D5","D8","D9","D10","D11","D16","D19","D25","D29","D33","D39","D41","D42","D44","D47","D48","D51")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.686.62"
$ws.Range("E2").Value = "  -1.55%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.596.19"
$ws.Range("E3").Value = "  -1.73%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "211.68"
$ws.Range("E5").Value = "  -1.43%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.56%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.12%  "

# Row 8 - was Cardano, now Dogecoin
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "0.0618"
$ws.Range("E8").Value = "  -2.19%  "

# Row 9 - was Dogecoin, now Cardano
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.247"
$ws.Range("E9").Value = "  -1.74%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.72"
$ws.Range("E10").Value = "  -2.01%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0836"
$ws.Range("E11").Value = "  -1.38%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.820.07"
$ws.Range("E12").Value = "  -1.70%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.609.00"
$ws.Range("E13").Value = "  -1.27%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -2.92%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -3.28%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "65.06"
$ws.Range("E16").Value = "  +0.44%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.658.48"
$ws.Range("E17").Value = "  -1.48%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -2.03%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "209.57"
$ws.Range("E19").Value = "  -2.25%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.18%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  -2.12%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -2.77%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  -2.30%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  -1.95%  "

# Row 25 - Monero
$ws.Range("D25").Value = "146.67"
$ws.Range("E25").Value = "  -0.96%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.24%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -2.58%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -0.36%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "15.34"
$ws.Range("E29").Value = "  -1.79%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -1.65%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.41%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -3.84%  "

# Row 33 - ImmutableX
$ws.Range("D33").Value = "0.671"
$ws.Range("E33").Value = "  -9.91%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -3.24%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.292.40"
$ws.Range("E35").Value = "  -5.01%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.78%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -5.97%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -3.50%  "

# Row 39 - ARBITRUM
$ws.Range("D39").Value = "0.835"
$ws.Range("E39").Value = "  -1.65%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  +0.11%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "0.793"
$ws.Range("E41").Value = "  -1.38%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "5.40"
$ws.Range("E42").Value = "  +0.90%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  -1.38%  "

# Row 44 - Aave
$ws.Range("D44").Value = "63.86"
$ws.Range("E44").Value = "  -1.90%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.732.35"
$ws.Range("E45").Value = "  -1.77%  "

# Row 46 - Quant
$ws.Range("E46").Value = "  -0.35%  "

# Row 47 - WEMIXToken
$ws.Range("D47").Value = "0.867"
$ws.Range("E47").Value = "  -1.44%  "

# Row 48 - RenderToken
$ws.Range("D48").Value = "1.62"
$ws.Range("E48").Value = "  -1.67%  "

# Row 49 - Algorand
$ws.Range("E49").Value = "  -3.29%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  -2.15%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "7.53"
$ws.Range("E51").Value = "  -2.16%  "
